$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sta")

# --- Step 1: Move existing "Note" text from column L to column M ---
$ws.Range("M10").Value = 'Multiple peaks over period range'
$ws.Range("M11").Value = 'Too few records'
$ws.Range("M12").Value = 'Multiple peaks over period range'
$ws.Range("M15").Value = 'Flat H/V and peak < 2'
$ws.Range("M17").Value = 'Multiple peaks over period range'
$ws.Range("M26").Value = 'Flat H/V and peak < 2'
$ws.Range("M29").Value = 'Too few records'
$ws.Range("M34").Value = 'Too few records'
$ws.Range("M37").Value = 'Flat H/V and peak < 2'
$ws.Range("M38").Value = 'Multiple peaks over period range'
$ws.Range("M40").Value = 'Multiple peaks over period range'
$ws.Range("M41").Value = 'Too few records'
$ws.Range("M61").Value = 'Flat H/V and peak < 2'
$ws.Range("M62").Value = 'Flat H/V and peak < 2'
$ws.Range("M63").Value = 'Flat H/V and peak < 2'

# --- Step 2: Set new numeric Reference Vs30 (USGS) values in column L ---
$ws.Range("L2").Value = 851
$ws.Range("L3").Value = 223
$ws.Range("L4").Value = 713
$ws.Range("L5").Value = 900
$ws.Range("L6").Value = 757
$ws.Range("L7").Value = 306
$ws.Range("L8").Value = 339
$ws.Range("L9").Value = 539
$ws.Range("L10").Value = 759
$ws.Range("L12").Value = 606
$ws.Range("L13").Value = 392
$ws.Range("L14").Value = 615
$ws.Range("L15").Value = 783
$ws.Range("L16").Value = 556
$ws.Range("L17").Value = 475
$ws.Range("L18").Value = 483
$ws.Range("L19").Value = 349
$ws.Range("L20").Value = 534
$ws.Range("L21").Value = 299
$ws.Range("L22").Value = 577
$ws.Range("L23").Value = 556
$ws.Range("L24").Value = 475
$ws.Range("L25").Value = 419
$ws.Range("L26").Value = 644
$ws.Range("L27").Value = 419
$ws.Range("L28").Value = 699
$ws.Range("L30").Value = 527
$ws.Range("L31").Value = 499
$ws.Range("L32").Value = 900
$ws.Range("L33").Value = 638
$ws.Range("L34").Value = 700
$ws.Range("L35").Value = 598
$ws.Range("L36").Value = 583
$ws.Range("L37").Value = 617
$ws.Range("L38").Value = 392
$ws.Range("L39").Value = 747
$ws.Range("L40").Value = 556
$ws.Range("L41").Value = 489
$ws.Range("L42").Value = 446
$ws.Range("L43").Value = 488
$ws.Range("L44").Value = 900
$ws.Range("L45").Value = 235
$ws.Range("L46").Value = 899
$ws.Range("L47").Value = 733
$ws.Range("L48").Value = 299
$ws.Range("L49").Value = 886
$ws.Range("L50").Value = 874
$ws.Range("L51").Value = 498
$ws.Range("L52").Value = 238
$ws.Range("L53").Value = 285
$ws.Range("L54").Value = 291
$ws.Range("L55").Value = 222
$ws.Range("L56").Value = 217
$ws.Range("L57").Value = 216
$ws.Range("L58").Value = 200
$ws.Range("L59").Value = 206
$ws.Range("L60").Value = 269
$ws.Range("L61").Value = 720
$ws.Range("L62").Value = 698
$ws.Range("L63").Value = 570
$ws.Range("L64").Value = 544
$ws.Range("L65").Value = 900
$ws.Range("L66").Value = 814
$ws.Range("L67").Value = 520
$ws.Range("L68").Value = 759

# --- Step 3: Clear column L for rows that only had a Note (no Vs30 data) ---
$ws.Range("L11").ClearContents()
$ws.Range("L29").ClearContents()

# --- Step 4: Update header row (L1 -> new ref Vs30 header, M1 -> "Note") ---
$ws.Range("M1").Value = "Note"
$ws.Range("L1").Value = "Reference Vs30 (USGS)"

# --- Step 5: Remove AutoFilter (and its cached sort state) ---
$ws.AutoFilterMode = $false

# --- Step 6: Update the hidden _FilterDatabase defined name to the new range ---
$fdb = $wb.Names.Item("sta!_FilterDatabase")
$fdb.RefersTo = '=sta!$A$1:$M$68'

# --- Step 7: Update the active selection ---
$ws.Range("L6").Select()

# --- Step 8: Reposition / resize the first picture (Picture 1) ---
$pic = $ws.Shapes.Item("Picture 1")
$pic.Left = 860.25
$pic.Top = 109.5
$pic.Width = 501.6873228346457
$pic.Height = 411.6985039370079

